# WhileLoops.pptx — "add while loops practice game"
#
# Slide 2 ("while Loops"), the code-block placeholder has a paragraph
#   <tab>// code
# The commit splits that run into:
#   <tab>            (same green/Consolas run, just shortened)
#   // body           (new green/Consolas run)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange

# Locate "// code" inside the full run of text for this placeholder.
$fullText = $tr.Text
$commentIdx0 = $fullText.IndexOf("// code")
if ($commentIdx0 -lt 0) {
    throw "Could not find '// code' in the Content Placeholder text"
}

# Characters() is 1-based; commentIdx0 (0-based) is already the 1-based
# start position of the "// code" substring.
$commentStart1 = $commentIdx0 + 1
$commentRange = $tr.Characters($commentStart1, 7)   # "// code" (7 chars)

# Replace "// code" with "// body" in place. Because this sub-range does
# not cover the leading tab character, the host splits the original run
# into the tab-only run (untouched formatting) and this new "// body"
# run, which inherits the same green/Consolas formatting.
$commentRange.Text = "// body"
